$wb = $excel.ActiveWorkbook
Write-Host $excel.CalculationVersion
try {
  $excel.CalculationVersion = 124519
  Write-Host "set ok"
} catch { Write-Host "ERR: $_" }
